# Actualización casos de uso
# Se eliminaron los casos de uso de las clases Tren, Colectivo y Subte
# (filas 17, 18 y 19 de la hoja: aplicarPrecio para Tren/Colectivo/Subte)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three rows corresponding to the Tren, Colectivo and Subte
# "aplicarPrecio" use cases (rows 17-19). This shifts all following rows
# up by three.
$ws.Rows("17:19").Delete()

# Restore the sequential "Nro. de Caso" numbering in column A for the
# rows that moved up.
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20

# Leave the selection where the editor ended up after the edit.
$ws.Range("C23").Select()
